$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1915.3334
$ws.Range("I43").Value = 1516.6666
$ws.Range("J43").Value = 2114.6667
$ws.Range("K43").Value = 1516.6666
$ws.Range("L43").Value = 2114.6667
$ws.Range("M43").Value = -1447.6666
$ws.Range("N43").Value = -2252.6667

$ws.Range("H51").Value = 4000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 4000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 4000
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = -4968

$ws.Range("H61").Value = 652
$ws.Range("I61").Value = 394.46155
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1183.38465
$ws.Range("L61").Value = 12000
$ws.Range("M61").Value = -1011.38465
$ws.Range("N61").Value = -12344

$ws.Range("H98").Value = 1324.95
$ws.Range("I98").Value = 1138.8334
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 1138.8334
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = 359.1666
$ws.Range("N98").Value = -5996

$ws.Range("H122").Value = 1324.95
$ws.Range("I122").Value = 1138.8334
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3416.5002
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -966.5001999999999
$ws.Range("N122").Value = -13900

$ws.Range("H123").Value = 34000
$ws.Range("J123").Value = 34000
$ws.Range("L123").Value = 34000
$ws.Range("N123").Value = -43800

$ws.Range("H125").Value = 1817.3334
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1817.3334
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 16356.0006
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = -21276.0006

$ws.Range("H127").Value = 1065.4117
$ws.Range("I127").Value = 441
$ws.Range("J127").Value = 1620.4445
$ws.Range("K127").Value = 1323
$ws.Range("L127").Value = 4861.333500000001
$ws.Range("M127").Value = 3637
$ws.Range("N127").Value = -14781.3335

$ws.Range("H132").Value = 3052.0857
$ws.Range("I132").Value = 2838.25
$ws.Range("J132").Value = 5333
$ws.Range("K132").Value = 8514.75
$ws.Range("L132").Value = 15999
$ws.Range("M132").Value = -5984.75
$ws.Range("N132").Value = -21059

$ws.Range("H137").Value = 1279.9818
$ws.Range("I137").Value = 1019.86487
$ws.Range("J137").Value = 1814.6666
$ws.Range("K137").Value = 3059.59461
$ws.Range("L137").Value = 5443.9998
$ws.Range("M137").Value = -509.5946100000001
$ws.Range("N137").Value = -10543.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 616.2
$ws.Range("I2").Value = 620
$ws.Range("K2").Value = 620
$ws.Range("M2").Value = -507

$ws.Range("H116").Value = 616.2
$ws.Range("I116").Value = 620
$ws.Range("K116").Value = 620
$ws.Range("M116").Value = 1674

$ws.Range("H122").Value = 16834634
$ws.Range("I122").Value = 201560
$ws.Range("K122").Value = 604680
$ws.Range("M122").Value = -602230

$ws.Range("H138").Value = 97585.60000000001
$ws.Range("J138").Value = 97585.60000000001
$ws.Range("L138").Value = 97585.60000000001
$ws.Range("N138").Value = -107865.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 616.2
$ws.Range("I3").Value = 620
$ws.Range("K3").Value = 620
$ws.Range("M3").Value = -506

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4943.398
$ws.Range("I31").Value = 1208.119
$ws.Range("J31").Value = 8353.869000000001
$ws.Range("K31").Value = 1208.119
$ws.Range("L31").Value = 8353.869000000001
$ws.Range("M31").Value = -913.1189999999999
$ws.Range("N31").Value = -8943.869000000001

$ws.Range("H34").Value = 4943.398
$ws.Range("I34").Value = 1208.119
$ws.Range("J34").Value = 8353.869000000001
$ws.Range("K34").Value = 1208.119
$ws.Range("L34").Value = 8353.869000000001
$ws.Range("M34").Value = -1006.119
$ws.Range("N34").Value = -8757.869000000001

$ws.Range("H44").Value = 6600
$ws.Range("J44").Value = 6600
$ws.Range("L44").Value = 6600
$ws.Range("N44").Value = -7484

$ws.Range("H54").Value = 6500
$ws.Range("I54").Value = 3000
$ws.Range("K54").Value = 3000
$ws.Range("M54").Value = -2342

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 7177.875
$ws.Range("J49").Value = 7177.875
$ws.Range("L49").Value = 21533.625
$ws.Range("N49").Value = -21845.625

$ws.Range("H54").Value = 11938.462
$ws.Range("J54").Value = 12766.667
$ws.Range("L54").Value = 38300.001
$ws.Range("N54").Value = -39418.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 23408.2
$ws.Range("I26").Value = 9499.75
$ws.Range("K26").Value = 9499.75
$ws.Range("M26").Value = -9219.75

$ws.Range("H50").Value = 23408.2
$ws.Range("I50").Value = 9499.75
$ws.Range("K50").Value = 9499.75
$ws.Range("M50").Value = -9001.75

$ws.Range("H141").Value = 70060.57000000001
$ws.Range("J141").Value = 70060.57000000001
$ws.Range("L141").Value = 70060.57000000001
$ws.Range("N141").Value = -80420.57000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1352.4117
$ws.Range("I46").Value = 927.2857
$ws.Range("J46").Value = 1650
$ws.Range("K46").Value = 927.2857
$ws.Range("L46").Value = 1650
$ws.Range("M46").Value = -739.2857
$ws.Range("N46").Value = -2026

$ws.Range("H93").Value = 10152.846
$ws.Range("I93").Value = 12228.3
$ws.Range("K93").Value = 12228.3
$ws.Range("M93").Value = -10980.3

$ws.Range("H132").Value = 3475.0688
$ws.Range("I132").Value = 3020.5264
$ws.Range("J132").Value = 4338.7
$ws.Range("K132").Value = 9061.5792
$ws.Range("L132").Value = 13016.1
$ws.Range("M132").Value = -6531.5792
$ws.Range("N132").Value = -18076.1

$ws.Range("H136").Value = 3705173.5
$ws.Range("I136").Value = 1282.579
$ws.Range("J136").Value = 23812010
$ws.Range("K136").Value = 3847.737
$ws.Range("L136").Value = 71436030
$ws.Range("M136").Value = -1297.737
$ws.Range("N136").Value = -71441130

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""

$ws.Range("H140").Value = 30057.375
$ws.Range("J140").Value = 30057.375
$ws.Range("L140").Value = 30057.375
$ws.Range("N140").Value = -40417.375

$ws.Range("H141").Value = 111588.08
$ws.Range("J141").Value = 111588.08
$ws.Range("L141").Value = 111588.08
$ws.Range("N141").Value = -121948.08
